# fix the bug of function SendMsg
# The FilePath column (F9:F14) on sheet "DataNode" pointed at
# ../resource/res/Scene/N.xml - the resource folder was renamed from
# "Scene" to "map", so the referenced paths need the same rename:
# ../resource/res/map/N.xml

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F9 -> plain replacement, no special formatting (matches original which
# also had no rich-text runs once fully replaced).
$ws.Range("F9").Value = "../resource/res/map/1.xml"

# F10:F14 keep the partially-rich-text look of the original cells (the
# "Scene"/"res" fragment had distinct run formatting) - retype the full
# path, then nudge the font on the whole cell (creates the new font/style)
# and again just on the "map" substring (creates the matching run in the
# shared string) so both the cell style and the rich-text run line up,
# the way Excel behaves when you edit part of a formatted cell.
$paths = @{
    10 = "../resource/res/map/2.xml"
    11 = "../resource/res/map/3.xml"
    12 = "../resource/res/map/4.xml"
    13 = "../resource/res/map/5.xml"
    14 = "../resource/res/map/6.xml"
}

foreach ($row in 10..14) {
    $cell = $ws.Cells.Item($row, 6)
    $text = $paths[$row]
    $cell.Value = $text
    $cell.Font.ColorIndex = 8
    $mapStart = $text.IndexOf("map") + 1
    $cell.Characters($mapStart, 3).Font.ColorIndex = 8
}

# Restore the selection recorded in the saved sheet view.
[void]$ws.Range("G18").Select()
